# Replay of the PW18-end-of-transcript edit: a reviewer comment is attached
# to the final "t" of "musket" in the caption "A 60 lb musket reaches 5
# <fr>pans</fr> and a half or 6", explaining that given the weight/size this
# is actually a "wall gun" rather than a handheld musket.
#
# The hyperlink r:id renumbering seen in the canonical XML diff (rId6->rId7,
# ..., rId13->rId14) is a mechanical side effect of a new relationship
# (word/comments.xml) being minted in the package -- it is not a content
# change we need to author explicitly; it falls out of adding the comment.

$d = $word.ActiveDocument

# Comments in this document are attributed to the reviewer, not "Word User".
$word.UserName = "Tillmann Taape"
$word.UserInitials = "TT"

# Locate the word "musket" in the caption text (unique in the document).
# Find.Execute mutates the range it is called against to the found span.
$musketRange = $d.Content
$found = $musketRange.Find.Execute("musket", $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find 'musket' in the document"
}

# The comment anchor only covers the final letter "t" of "musket" (i.e. the
# run "musket" is split into "muske" + "t", and the commentRange wraps just
# the trailing "t").
$commentAnchor = $d.Range($musketRange.End - 1, $musketRange.End)

$commentText = "JT: given the weight and size, this is not a handheld musket, but what was known as a wall gun."

$comment = $d.Comments.Add($commentAnchor, $commentText)
